$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "model_2_4_0"
$ws.Range("B2").Value = -0.001236788884735551
$ws.Range("C2").Value = -0.0771084754329987
$ws.Range("D2").Value = -0.341152233302648
$ws.Range("E2").Value = -0.05115300178463333
$ws.Range("F2").Value = 1.108074069023132
$ws.Range("G2").Value = 1.959989070892334
$ws.Range("H2").Value = 0.2471889555454254
$ws.Range("I2").Value = 1.15396523475647

$ws.Range("A3").Value = "model_2_4_1"
$ws.Range("B3").Value = 0.03860672347268723
$ws.Range("C3").Value = -0.1027274069571471
$ws.Range("D3").Value = -0.3719719032556219
$ws.Range("E3").Value = -0.07606939401569046
$ws.Range("F3").Value = 1.063979029655457
$ws.Range("G3").Value = 2.006607294082642
$ws.Range("H3").Value = 0.2528693675994873
$ws.Range("I3").Value = 1.181318759918213

$ws.Range("A4").Value = "model_2_4_14"
$ws.Range("B4").Value = 0.09919338355202301
$ws.Range("C4").Value = -0.1967812923922139
$ws.Range("D4").Value = -0.4117062867951791
$ws.Range("E4").Value = -0.161743799631052
$ws.Range("F4").Value = 0.9969275593757629
$ws.Range("G4").Value = 2.177754878997803
$ws.Range("H4").Value = 0.2601928412914276
$ws.Range("I4").Value = 1.275372743606567

$ws.Range("A5").Value = "model_2_4_2"
$ws.Range("B5").Value = 0.09929134637632209
$ws.Range("C5").Value = -0.1467233699917059
$ws.Range("D5").Value = -0.5171858507911891
$ws.Range("E5").Value = -0.1261501727855825
$ws.Range("F5").Value = 0.9968190789222717
$ws.Range("G5").Value = 2.086665630340576
$ws.Range("H5").Value = 0.2796338796615601
$ws.Range("I5").Value = 1.236297845840454

$ws.Range("A6").Value = "model_2_4_16"
$ws.Range("B6").Value = 0.09930816734455727
$ws.Range("C6").Value = -0.196515009631226
$ws.Range("D6").Value = -0.4099056673005024
$ws.Range("E6").Value = -0.1613678037111921
$ws.Range("F6").Value = 0.9968004822731018
$ws.Range("G6").Value = 2.177270412445068
$ws.Range("H6").Value = 0.2598609328269958
$ws.Range("I6").Value = 1.274960160255432

$ws.Range("A7").Value = "model_2_4_15"
$ws.Range("B7").Value = 0.0993141001281147
$ws.Range("C7").Value = -0.196500620812958
$ws.Range("D7").Value = -0.4098136995747119
$ws.Range("E7").Value = -0.1613478982550252
$ws.Range("F7").Value = 0.9967938661575317
$ws.Range("G7").Value = 2.177244186401367
$ws.Range("H7").Value = 0.2598440051078796
$ws.Range("I7").Value = 1.274938106536865

$ws.Range("A8").Value = "model_2_4_22"
$ws.Range("B8").Value = 0.09932080507882668
$ws.Range("C8").Value = -0.1964480985326931
$ws.Range("D8").Value = -0.4101482731144555
$ws.Range("E8").Value = -0.1613282421129689
$ws.Range("F8").Value = 0.9967864155769348
$ws.Range("G8").Value = 2.177148580551147
$ws.Range("H8").Value = 0.2599056661128998
$ws.Range("I8").Value = 1.274916648864746

$ws.Range("A9").Value = "model_2_4_21"
$ws.Range("B9").Value = 0.09932080507882668
$ws.Range("C9").Value = -0.1964480985326931
$ws.Range("D9").Value = -0.4101482731144555
$ws.Range("E9").Value = -0.1613282421129689
$ws.Range("F9").Value = 0.9967864155769348
$ws.Range("G9").Value = 2.177148580551147
$ws.Range("H9").Value = 0.2599056661128998
$ws.Range("I9").Value = 1.274916648864746

$ws.Range("A10").Value = "model_2_4_20"
$ws.Range("B10").Value = 0.09932080507882668
$ws.Range("C10").Value = -0.1964480985326931
$ws.Range("D10").Value = -0.4101482731144555
$ws.Range("E10").Value = -0.1613282421129689
$ws.Range("F10").Value = 0.9967864155769348
$ws.Range("G10").Value = 2.177148580551147
$ws.Range("H10").Value = 0.2599056661128998
$ws.Range("I10").Value = 1.274916648864746

$ws.Range("A11").Value = "model_2_4_19"
$ws.Range("B11").Value = 0.09932080507882668
$ws.Range("C11").Value = -0.1964480985326931
$ws.Range("D11").Value = -0.4101482731144555
$ws.Range("E11").Value = -0.1613282421129689
$ws.Range("F11").Value = 0.9967864155769348
$ws.Range("G11").Value = 2.177148580551147
$ws.Range("H11").Value = 0.2599056661128998
$ws.Range("I11").Value = 1.274916648864746

$ws.Range("A12").Value = "model_2_4_18"
$ws.Range("B12").Value = 0.09932080507882668
$ws.Range("C12").Value = -0.1964480985326931
$ws.Range("D12").Value = -0.4101482731144555
$ws.Range("E12").Value = -0.1613282421129689
$ws.Range("F12").Value = 0.9967864155769348
$ws.Range("G12").Value = 2.177148580551147
$ws.Range("H12").Value = 0.2599056661128998
$ws.Range("I12").Value = 1.274916648864746

$ws.Range("A13").Value = "model_2_4_17"
$ws.Range("B13").Value = 0.09932080507882668
$ws.Range("C13").Value = -0.1964480985326931
$ws.Range("D13").Value = -0.4101482731144555
$ws.Range("E13").Value = -0.1613282421129689
$ws.Range("F13").Value = 0.9967864155769348
$ws.Range("G13").Value = 2.177148580551147
$ws.Range("H13").Value = 0.2599056661128998
$ws.Range("I13").Value = 1.274916648864746

$ws.Range("A14").Value = "model_2_4_24"
$ws.Range("B14").Value = 0.09932080507882668
$ws.Range("C14").Value = -0.1964480985326931
$ws.Range("D14").Value = -0.4101482731144555
$ws.Range("E14").Value = -0.1613282421129689
$ws.Range("F14").Value = 0.9967864155769348
$ws.Range("G14").Value = 2.177148580551147
$ws.Range("H14").Value = 0.2599056661128998
$ws.Range("I14").Value = 1.274916648864746

$ws.Range("A15").Value = "model_2_4_23"
$ws.Range("B15").Value = 0.09932080507882668
$ws.Range("C15").Value = -0.1964480985326931
$ws.Range("D15").Value = -0.4101482731144555
$ws.Range("E15").Value = -0.1613282421129689
$ws.Range("F15").Value = 0.9967864155769348
$ws.Range("G15").Value = 2.177148580551147
$ws.Range("H15").Value = 0.2599056661128998
$ws.Range("I15").Value = 1.274916648864746

$ws.Range("A16").Value = "model_2_4_12"
$ws.Range("B16").Value = 0.09939638384743188
$ws.Range("C16").Value = -0.195432282902307
$ws.Range("D16").Value = -0.4179740990600016
$ws.Range("E16").Value = -0.161055069300692
$ws.Range("F16").Value = 0.9967028498649597
$ws.Range("G16").Value = 2.175300121307373
$ws.Range("H16").Value = 0.2613480687141418
$ws.Range("I16").Value = 1.274616599082947

$ws.Range("A17").Value = "model_2_4_10"
$ws.Range("B17").Value = 0.0994318744978907
$ws.Range("C17").Value = -0.193277522791294
$ws.Range("D17").Value = -0.4384853791436401
$ws.Range("E17").Value = -0.1607848138527508
$ws.Range("F17").Value = 0.9966636300086975
$ws.Range("G17").Value = 2.171379089355469
$ws.Range("H17").Value = 0.265128493309021
$ws.Range("I17").Value = 1.274319887161255

$ws.Range("A18").Value = "model_2_4_13"
$ws.Range("B18").Value = 0.09962912926329348
$ws.Range("C18").Value = -0.1952749726472431
$ws.Range("D18").Value = -0.4100658353433144
$ws.Range("E18").Value = -0.1602922780579776
$ws.Range("F18").Value = 0.9964452385902405
$ws.Range("G18").Value = 2.175013780593872
$ws.Range("H18").Value = 0.2598904669284821
$ws.Range("I18").Value = 1.273779273033142

$ws.Range("A19").Value = "model_2_4_11"
$ws.Range("B19").Value = 0.09981486051873822
$ws.Range("C19").Value = -0.1929980713590784
$ws.Range("D19").Value = -0.426604944680764
$ws.Range("E19").Value = -0.1596009409579628
$ws.Range("F19").Value = 0.9962397813796997
$ws.Range("G19").Value = 2.170870780944824
$ws.Range("H19").Value = 0.2629387974739075
$ws.Range("I19").Value = 1.273020386695862

$ws.Range("A20").Value = "model_2_4_9"
$ws.Range("B20").Value = 0.100859368823372
$ws.Range("C20").Value = -0.1857795595441754
$ws.Range("D20").Value = -0.4542717490203676
$ws.Range("E20").Value = -0.1554524533141999
$ws.Range("F20").Value = 0.995083749294281
$ws.Range("G20").Value = 2.157735347747803
$ws.Range("H20").Value = 0.2680381238460541
$ws.Range("I20").Value = 1.268465995788574

$ws.Range("A21").Value = "model_2_4_8"
$ws.Range("B21").Value = 0.1014686107787831
$ws.Range("C21").Value = -0.1789144520370398
$ws.Range("D21").Value = -0.4855427840400852
$ws.Range("E21").Value = -0.1518987257648128
$ws.Range("F21").Value = 0.994409441947937
$ws.Range("G21").Value = 2.145242929458618
$ws.Range("H21").Value = 0.2738017141819
$ws.Range("I21").Value = 1.264564752578735

$ws.Range("A22").Value = "model_2_4_3"
$ws.Range("B22").Value = 0.1020292019527212
$ws.Range("C22").Value = -0.1702853267273956
$ws.Range("D22").Value = -0.5120716487583978
$ws.Range("E22").Value = -0.1464224385256563
$ws.Range("F22").Value = 0.9937890768051147
$ws.Range("G22").Value = 2.129540681838989
$ws.Range("H22").Value = 0.2786912620067596
$ws.Range("I22").Value = 1.25855278968811

$ws.Range("A23").Value = "model_2_4_4"
$ws.Range("B23").Value = 0.1021750504245095
$ws.Range("C23").Value = -0.1701528952795848
$ws.Range("D23").Value = -0.5128949163892607
$ws.Range("E23").Value = -0.146371268476196
$ws.Range("F23").Value = 0.993627667427063
$ws.Range("G23").Value = 2.129300117492676
$ws.Range("H23").Value = 0.2788430154323578
$ws.Range("I23").Value = 1.258496642112732

$ws.Range("A24").Value = "model_2_4_5"
$ws.Range("B24").Value = 0.1026596559564265
$ws.Range("C24").Value = -0.1718205537665856
$ws.Range("D24").Value = -0.4865137808111959
$ws.Range("E24").Value = -0.1457503463690848
$ws.Range("F24").Value = 0.9930914044380188
$ws.Range("G24").Value = 2.132334470748901
$ws.Range("H24").Value = 0.2739806473255157
$ws.Range("I24").Value = 1.257815003395081

$ws.Range("A25").Value = "model_2_4_6"
$ws.Range("B25").Value = 0.1027171165504086
$ws.Range("C25").Value = -0.1724941609164194
$ws.Range("D25").Value = -0.4811592322862823
$ws.Range("E25").Value = -0.1459184729065193
$ws.Range("F25").Value = 0.9930276870727539
$ws.Range("G25").Value = 2.133560419082642
$ws.Range("H25").Value = 0.2729937732219696
$ws.Range("I25").Value = 1.257999539375305

$ws.Range("A26").Value = "model_2_4_7"
$ws.Range("B26").Value = 0.1032754513636238
$ws.Range("C26").Value = -0.172439474490192
$ws.Range("D26").Value = -0.4634053168027545
$ws.Range("E26").Value = -0.1444677694639867
$ws.Range("F26").Value = 0.9924098253250122
$ws.Range("G26").Value = 2.133460521697998
$ws.Range("H26").Value = 0.2697215378284454
$ws.Range("I26").Value = 1.256407141685486
